$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 14: Qty changed from 1.5 to 3 (dependent formulas recalc automatically)
$ws.Range("B14").Value = 3

# Row 45: item reference text update
$ws.Range("C45").Value = "EA77-625"

# Row 45: Value (G45) updated, Companie (H45) and Who? (M45) filled in
$ws.Range("G45").Value = 47600
$ws.Range("H45").Value = "Sigma"
$ws.Range("M45").Value = "Johan"

# Update the saved view state (scroll position / active selection)
$ws.Range("G45").Select()
$wv = $excel.ActiveWindow
$wv.ScrollRow = 27
$wv.ScrollColumn = 1

$wb.Save()
